{"js": "// Update the video link document:\n//  - youtu.be short link gets a new video id + the share \"?si=\" query string\n//  - the blank paragraph between the short link and the instructions is removed\n//  - the \"url hit in the web broweser to open the video \" paragraph is\n//    normalized to a single run (the spell-check proofErr markers / run\n//    splits around \"url\" and \"broweser\" disappear)\n//  - the youtube.com/watch link gets the new video id\n\nconst body = context.document.body;\n\n// Step 1: update the youtu.be short link, adding the new share query string\nconst shortResults = body.search(\"youtu.be/2yYx-dsY5Us\", { matchCase: true });\nshortResults.load(\"items\");\nawait context.sync();\nif (shortResults.items.length > 0) {\n  shortResults.items[0].insertText(\n    \"youtu.be/T9xb_Vy9_UY?si=fj6_uatCEy6h0V-2\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// Step 2: remove the stray empty paragraph right after the short link\nconst paras1 = body.paragraphs;\nparas1.load(\"items/text\");\nawait context.sync();\nlet emptyPara = null;\nfor (const p of paras1.items) {\n  if (p.text.trim() === \"\") {\n    emptyPara = p;\n    break;\n  }\n}\nif (emptyPara) {\n  emptyPara.delete();\n  await context.sync();\n}\n\n// Step 3: collapse the \"url hit in the web broweser...\" paragraph's runs\n// (and the proofErr spell-check markers around \"url\"/\"broweser\") into a\n// single clean run with the same text. Inserting a fresh paragraph with the\n// plain text and deleting the old one drops the now-orphaned proofErr marks.\nconst paras2 = body.paragraphs;\nparas2.load(\"items/text\");\nawait context.sync();\nlet instrPara = null;\nfor (const p of paras2.items) {\n  if (p.text.indexOf(\"hit in the web broweser\") !== -1) {\n    instrPara = p;\n    break;\n  }\n}\nif (instrPara) {\n  instrPara.insertParagraph(\n    \"url hit in the web broweser to open the video \",\n    Word.InsertLocation.after\n  );\n  instrPara.delete();\n  await context.sync();\n}\n\n// Step 4: update the youtube.com watch link with the new video id\nconst longResults = body.search(\"watch?v=2yYx-dsY5Us\", { matchCase: true });\nlongResults.load(\"items\");\nawait context.sync();\nif (longResults.items.length > 0) {\n  longResults.items[0].insertText(\"watch?v=T9xb_Vy9_UY\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the video link document:\n#  - youtu.be short link gets a new video id + the share \"?si=\" query string\n#  - the blank paragraph between the short link and the instructions is removed\n#  - the \"url hit in the web broweser to open the video \" paragraph is\n#    normalized to a single run (the spell-check proofErr markers / run\n#    splits around \"url\" and \"broweser\" disappear)\n#  - the youtube.com/watch link gets the new video id\n\n$d = $word.ActiveDocument\n\n$oldId = \"2yYx-dsY5Us\"\n$newId = \"T9xb_Vy9_UY\"\n$newQuery = \"?si=fj6_uatCEy6h0V-2\"\n\n# Step 1: update the youtu.be short-link paragraph, appending the new share query string\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*youtu.be/$oldId*\") {\n        $p.Range.Find.Execute($oldId, $false, $false, $false, $false, $false, $true, 1, $false, \"$newId$newQuery\", 2) | Out-Null\n        break\n    }\n}\n\n# Step 2: remove the stray empty paragraph right after the short link\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# Step 3: collapse the \"url hit in the web broweser...\" paragraph's runs\n# (and the proofErr spell-check markers around \"url\"/\"broweser\") into a\n# single clean run with the same text\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*hit in the web broweser*\") {\n        $full = $p.Range.Text.TrimEnd([char]13)\n        $start = $p.Range.Start\n        $endInclMark = $p.Range.End\n        $d.Range($start, $endInclMark).Delete()\n        $d.Range($start, $start).InsertAfter($full + \"`r\")\n        break\n    }\n}\n\n# Step 4: update the full youtube.com watch-link paragraph with the new video id\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*watch?v=$oldId*\") {\n        $p.Range.Find.Execute($oldId, $false, $false, $false, $false, $false, $true, 1, $false, $newId, 2) | Out-Null\n        break\n    }\n}\n"}
